$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PATRIKA 25-26")

$ws.Range("A4").Value = "1-Jul-25 to 21-Jan-26"
$ws.Range("B7").Value = "1-Jul-25 to 21-Jan-26"
$ws.Range("A636").Value = "9333 CARDS-FRIDAY (6042)"
$ws.Range("B35").Value = 31
$ws.Range("C35").Value = 108.5
$ws.Range("E35").Value = 249.55
$ws.Range("B44").Value = 23
$ws.Range("C44").Value = 49
$ws.Range("E44").Value = 130.34
$ws.Range("B144").Value = 14
$ws.Range("C144").Value = 1
$ws.Range("E144").Value = 18.5
$ws.Range("B162").Value = 13
$ws.Range("C162").Value = 1.49
$ws.Range("E162").Value = 17.14
$ws.Range("B195").Value = 32
$ws.Range("C195").Value = 68.5
$ws.Range("E195").Value = 308.25
$ws.Range("B200").Value = 17
$ws.Range("C200").Value = 98.5
$ws.Range("E200").Value = 334.9
$ws.Range("B207").Value = 45
$ws.Range("C207").Value = 55.5
$ws.Range("E207").Value = 217.67
$ws.Range("B241").Value = 81
$ws.Range("C241").Value = 13.5
$ws.Range("E241").Value = 85.05
$ws.Range("B249").Value = 3
$ws.Range("C249").Value = 22.5
$ws.Range("E249").Value = 135
$ws.Range("C250").Value = 107.5
$ws.Range("E250").Value = 376.25
$ws.Range("B274").Value = 20
$ws.Range("C274").Value = 26.5
$ws.Range("E274").Value = 172.25
$ws.Range("B275").Value = 10
$ws.Range("C275").Value = 23.5
$ws.Range("E275").Value = 122.91
$ws.Range("B289").Value = 49
$ws.Range("C289").Value = 179
$ws.Range("E289").Value = 874.97
$ws.Range("B309").Value = 27
$ws.Range("C309").Value = 33.450000000000003
$ws.Range("E309").Value = 316.10000000000002
$ws.Range("B325").Value = 36
$ws.Range("C325").Value = 170.91
$ws.Range("E325").Value = 731.49
$ws.Range("B337").Value = 20
$ws.Range("C337").Value = 4
$ws.Range("E337").Value = 36
$ws.Range("B360").Value = 33
$ws.Range("C360").Value = 36
$ws.Range("E360").Value = 213.84
$ws.Range("B369").Value = 41
$ws.Range("C369").Value = 46
$ws.Range("E369").Value = 377.2
$ws.Range("B370").Value = 112
$ws.Range("C370").Value = 16.45
$ws.Range("E370").Value = 74.03
$ws.Range("B372").Value = 5
$ws.Range("C372").Value = 20.5
$ws.Range("E372").Value = 188.6
$ws.Range("B374").Value = 11
$ws.Range("C374").Value = 45
$ws.Range("E374").Value = 315
$ws.Range("B385").Value = 6
$ws.Range("C385").Value = 20.5
$ws.Range("E385").Value = 116.85
$ws.Range("B408").Value = 61
$ws.Range("C408").Value = 44.5
$ws.Range("E408").Value = 231.4
$ws.Range("C445").Value = 12.5
$ws.Range("E445").Value = 131.25
$ws.Range("C461").Value = 30.16
$ws.Range("E461").Value = 257.87
$ws.Range("C470").Value = 31.5
$ws.Range("E470").Value = 299.25
$ws.Range("B475").Value = 39
$ws.Range("C475").Value = 23
$ws.Range("E475").Value = 237.39
$ws.Range("B476").Value = 6
$ws.Range("C476").Value = 4.5
$ws.Range("E476").Value = 58.5
$ws.Range("B487").Value = 5
$ws.Range("C487").Value = 5.5
$ws.Range("E487").Value = 55
$ws.Range("B537").Value = 19
$ws.Range("C537").Value = 31.5
$ws.Range("E537").Value = 63
$ws.Range("B579").Value = 57
$ws.Range("C579").Value = 287.5
$ws.Range("E579").Value = 431.25
$ws.Range("B580").Value = 89
$ws.Range("C580").Value = 251
$ws.Range("B581").Value = 60
$ws.Range("C581").Value = 315.5
$ws.Range("B582").Value = 80
$ws.Range("C582").Value = 102.5
$ws.Range("B589").Value = 93
$ws.Range("C589").Value = 933
$ws.Range("E589").Value = 793.05
$ws.Range("B607").Value = 14
$ws.Range("C607").Value = 9
$ws.Range("E607").Value = 90
$ws.Range("B623").Value = 32
$ws.Range("C623").Value = 22.5
$ws.Range("E623").Value = 33.75
$ws.Range("B626").Value = 62
$ws.Range("C626").Value = 107.5
$ws.Range("E626").Value = 245.1
$ws.Range("B636").Value = 9
$ws.Range("C636").Value = -0.5
$ws.Range("D636").Value = 5.46
$ws.Range("E636").Value = -2.73
$ws.Range("B638").Value = 95
$ws.Range("C638").Value = 263.89
$ws.Range("E638").Value = 448.61
$ws.Range("B650").Value = 28
$ws.Range("C650").Value = 10
$ws.Range("E650").Value = 31
$ws.Range("B660").Value = 123
$ws.Range("C660").Value = 19.5
$ws.Range("E660").Value = 60.45
$ws.Range("B661").Value = 39
$ws.Range("C661").Value = 1
$ws.Range("E661").Value = 3.6
$ws.Range("B664").Value = 36
$ws.Range("C664").Value = 21
$ws.Range("E664").Value = 78.75
$ws.Range("C670").Value = 145
$ws.Range("E670").Value = 326.25
$ws.Range("C680").Value = 7.5
$ws.Range("E680").Value = 37.5
$ws.Range("B684").Value = 12
$ws.Range("C684").Value = 12
$ws.Range("B718").Value = 47
$ws.Range("C718").Value = 68.75
$ws.Range("E718").Value = 171.88
$ws.Range("B719").Value = 58
$ws.Range("C719").Value = 56.75
$ws.Range("E719").Value = 141.88
$ws.Range("C721").Value = 43767.28
$ws.Range("E721").Value = 106398.51
